$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (G1, H1) with the same header style as A1:F1 ---
$ws.Range("A1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("G1").Value = "Elapsed Time"
$ws.Range("H1").Value = "CPU"

# --- Updated metric values in row 2 ---
$ws.Range("B2").Value = 0.7169911074822892
$ws.Range("C2").Value = 0.9789400929690297
$ws.Range("D2").Value = 0.6525923309357241

# --- Updated model description text ---
$ws.Range("F2").Value = "Pipeline(steps=[('model', AdaBoostRegressor(learning_rate=0.5))])"

# --- New elapsed time / CPU values ---
$ws.Range("G2").Value = 0.1228190763666741
$ws.Range("H2").Value = 0.991
